# Add newly captured Evernote links to the hyperlinks tracking sheet.
# Column A = note title, Column B = the Evernote share link (displayed as
# the raw URL, same text as the hyperlink target), styled with the
# workbook's "Hyperlink" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-LinkRow {
    param([int]$RowNumber, [string]$Title, [string]$Url)

    $ws.Range("A$RowNumber").Value = $Title

    $cell = $ws.Range("B$RowNumber")
    $null = $ws.Hyperlinks.Add($cell, $Url, "", "", $Url)
    $cell.Style = "Hyperlink"
}

# Row 675 is entered first (it lands at the bottom of the sheet but was the
# first new hyperlink relationship created).
Add-LinkRow 675 "[ESG] Dialoca - Gente e Gestão" "https://www.evernote.com/shard/s567/nl/98932539/46cd58c0-f114-451e-b88a-4d117566c14a?title=%5BESG%5D%20Dialoca%20-%20Gente%20e%20Gest%C3%A3o"

Add-LinkRow 668 "Bluefit | Conversa com Franqueador" "https://www.evernote.com/shard/s567/nl/98932539/f8400aef-8aee-4561-8e24-a57b9b5097b4?title=Bluefit%20%7C%20Conversa%20com%20Franqueador"

Add-LinkRow 669 "Wine | Round Table" "https://www.evernote.com/shard/s567/nl/98932539/6e26f090-32a1-4b5e-c922-19faf3622489?title=Wine%20%7C%20Round%20Table"

Add-LinkRow 670 "CVC | Call com RI" "https://www.evernote.com/shard/s567/nl/98932539/be4a2898-aab9-ce4d-92a4-62de1386499a?title=CVC%20%7C%20Call%20com%20RI"

Add-LinkRow 671 "CVC | Papo com Pessoa do Setor" "https://www.evernote.com/shard/s567/nl/98932539/5569ef6a-d2d1-56b3-b15e-dd5e210c7b25?title=CVC%20%7C%20Papo%20com%20Pessoa%20do%20Setor"

Add-LinkRow 672 "Selfit | Conversa com CFO sobre a empresa e mercado de academias" "https://www.evernote.com/shard/s567/nl/98932539/0f297589-cbe2-410b-a1dd-c37c73c95084?title=Selfit%20%7C%20Conversa%20com%20CFO%20sobre%20a%20empresa%20e%20mercado%20de%20academias"

Add-LinkRow 673 "Heverton Peixoto | Wiz CEO" "https://www.evernote.com/shard/s567/nl/98932539/ecf04a2f-7f94-4e5a-9e8e-b82d9660b091?title=Heverton%20Peixoto%20%7C%20Wiz%20CEO"

Add-LinkRow 674 "Smart Fit | IE com Santander" "https://www.evernote.com/shard/s567/nl/98932539/444fc6f6-cb6c-4d58-ade7-ccf7fe40be72?title=Smart%20Fit%20%7C%20IE%20com%20Santander"

# Move the selection/scroll to where the user ended up after pasting the
# new rows at the bottom of the sheet.
$ws.Range("A675").Select()
